$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 30 into a new row, shifting the existing row 30 (and everything
# below it) down by one. This mirrors the underlying edit: a new weekly
# record is inserted before the old row 30, pushing rows 30:58 to 31:59.
$ws.Rows.Item(30).Copy()
$ws.Rows.Item(30).Insert()

# The freshly inserted row 30 is an exact duplicate of (what is now) row 31.
# Update the two fields that actually differ for the new record: the date
# and the origin.
$ws.Range("D30").Value = 44580
$ws.Range("O30").Value = "Provincia del Elquí"
